$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row before the existing row 127, shifting rows 127-198
# down to 128-199 (dimension grows from A1:R198 to A1:R199).
$ws.Rows("127:127").Insert()

$ws.Range("A127").Value = 6
$ws.Range("B127").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C127").Value = "Metropolitana"
$ws.Range("D127").Value = 44603
$ws.Range("E127").Value = 13
$ws.Range("F127").Value = 100112022
$ws.Range("G127").Value = "Arveja Verde"
$ws.Range("H127").Value = "Sin especificar"
$ws.Range("I127").Value = "Primera"
$ws.Range("J127").Value = 200
$ws.Range("K127").Value = 25000
$ws.Range("L127").Value = 26000
$ws.Range("M127").Value = 25600
$ws.Range("N127").Value = "$/saco 25 kilos"
$ws.Range("O127").Value = "Carahue"
$ws.Range("P127").Value = 1024
$ws.Range("Q127").Value = 25
$ws.Range("R127").Value = "Hortaliza"
